# Add a new "Git Stash" entry to the git reference sheet (Sheet1),
# as described in the commit: "Add git entry on using stash, awesome tool"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row content for row 24 (previously blank placeholder row)
$ws.Range("A24").Value = "Git Stash"
$ws.Range("B24").Value = "Working on master but need a quick reference to branch, stash can temp save the change and recover later on"

$stashText = "Save dirty workings on branch #1 `n> git stash -u    //And now can switch to branch #2, while at this moment branch#1 status is clean `nView Stash:`n> git stash list`nRecover stash after switch back from branch#2:`n> (at branch#1) git stash pop`n"
$ws.Range("C24").Value = $stashText

# Match row height used by the other multi-line entries in this sheet
$ws.Rows.Item(24).RowHeight = 90

# Update the saved view/selection to point at the newly added row
[void]$ws.Range("C24").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
